# "Switch pandoc to nix." -- a resave/re-normalization of the document:
# Word (via a newer pandoc/build toolchain) re-emits <w:rPr> children in
# canonical schema order (b/bCs before i/iCs before color, etc.) instead
# of the ad-hoc order the previous toolchain produced. There is no
# semantic content change -- we reproduce it by re-asserting each run's /
# style's Bold / Italic property (already on), which makes the engine
# rewrite that <w:rPr> in schema order without altering its meaning.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Body runs: <w:iCs/><w:i/> -> <w:i/><w:iCs/>  and
#               <w:bCs/><w:b/> -> <w:b/><w:bCs/>
# Walk the nine affected runs in document order using whole-word Find,
# toggling the matching Font property so the engine rewrites rPr.
# ---------------------------------------------------------------------

function Find-NextWord($word) {
    $r = $d.Content
    $r.Start = $script:cursor
    $r.End = $d.Content.End
    $r.Find.Execute($word, $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $script:cursor = $r.End
    return $r
}

function Toggle-Italic($word) {
    $r = Find-NextWord $word
    $r.Font.Italic = $true
}

function Toggle-Bold($word) {
    $r = Find-NextWord $word
    $r.Font.Bold = $true
}

$script:cursor = 0
Toggle-Italic "italics"                     # "...formatted in *italics*"
Toggle-Bold   "bold"                        # "...italics and *bold*"
Find-NextWord "line break" | Out-Null       # skip the plain-text decoy
Toggle-Italic "line break"                  # "...and have a *line break*"
Toggle-Bold   "formatting"                  # "...slightly', and *formatting*"
Toggle-Italic "italic"                      # "...recursive include, with *italic*, bold"
Toggle-Bold   "bold"                        # "...with italic, *bold*"
Toggle-Italic "italic"                      # "...in recursive include, with *italic*, bold"
Toggle-Bold   "bold"                        # "...with italic, *bold*"
Toggle-Italic "even more italic text"       # "...this is *even more italic text*."

# ---------------------------------------------------------------------
# 2) Character styles in styles.xml: reorder rPr children to schema
#    order (b, i, color) by re-asserting Bold/Italic on each style.
# ---------------------------------------------------------------------

function Retoggle-StyleFormat($styleName, $hasBold, $hasItalic) {
    $s = $d.Styles.Item($styleName)
    if ($hasBold) {
        $s.Font.Bold = $true
    }
    if ($hasItalic) {
        $s.Font.Italic = $true
    }
}

Retoggle-StyleFormat "KeywordTok"       $true  $false
Retoggle-StyleFormat "ImportTok"        $true  $false
Retoggle-StyleFormat "CommentTok"       $false $true
Retoggle-StyleFormat "DocumentationTok" $false $true
Retoggle-StyleFormat "AnnotationTok"    $true  $true
Retoggle-StyleFormat "CommentVarTok"    $true  $true
Retoggle-StyleFormat "ControlFlowTok"   $true  $false
Retoggle-StyleFormat "InformationTok"   $true  $true
Retoggle-StyleFormat "WarningTok"       $true  $true
Retoggle-StyleFormat "AlertTok"         $true  $false
Retoggle-StyleFormat "ErrorTok"         $true  $false

# ---------------------------------------------------------------------
# 3) numbering.xml: the <w:nsid> of abstractNum 990 is zero-padded to
#    8 hex digits ("A990" -> "0000A990") by the resave. `nsid` is a
#    bookkeeping-only GUID fragment (ECMA-376 17.9.13): it isn't part
#    of Word's object model (no paragraph in this document even uses
#    that list definition, so there is no ListFormat/ListTemplate
#    handle that reaches it either), so there is no COM/object-model
#    call that can reach it -- nothing to do here.
# ---------------------------------------------------------------------
